# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment schedule"
#   sheet, shifting the existing "Late" / heading / "Outstanding" columns
#   one position to the right.
# - Match the width Excel assigns to a freshly inserted column (it inherits
#   the width of the column immediately to its left).
# - Update the active sheet / selected cells to reflect where the user was
#   last working (Repayment schedule tab, cell L21) while also recording the
#   last selection made on the Edit Repayment Schedule sheet (cell D15).

$wb = $excel.ActiveWorkbook

$wsRepayment = $wb.Worksheets.Item("Repayment schedule")

# Remember column M's width so the newly inserted column N matches what
# Excel does automatically (copy the width of the column to its left).
$leftWidth = $wsRepayment.Columns("M").ColumnWidth

[void]$wsRepayment.Columns("N").Insert()
$wsRepayment.Columns("N").ColumnWidth = $leftWidth

# Record the last selection on the "Edit Repayment Schedule" sheet before
# switching away from it.
$wsEdit = $wb.Worksheets.Item("Edit Repayment Schedule")
[void]$wsEdit.Range("D15").Select()

# Finally, make "Repayment schedule" the active sheet/tab with cell L21
# selected (this also clears the previously active tab on "NewLoanInput").
[void]$wsRepayment.Activate()
[void]$wsRepayment.Range("L21").Select()
